$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-19 Sunday" "2025-10-20 Monday"

Replace-Text "445×8=3560" "637×2=1274"
Replace-Text "238×8=1904" "947×3=2841"
Replace-Text "880×3=2640" "266×3=798"
Replace-Text "901×5=4505" "203×3=609"
Replace-Text "297×3=891" "710×6=4260"

Replace-Text "507×7=3549" "375×7=2625"
Replace-Text "842×6=5052" "159×7=1113"
Replace-Text "694×7=4858" "145×5=725"
Replace-Text "289×7=2023" "928×9=8352"
Replace-Text "331×9=2979" "891×8=7128"

Replace-Text "804×3=2412" "141×5=705"
Replace-Text "287×3=861" "510×4=2040"
Replace-Text "569×9=5121" "546×5=2730"
Replace-Text "955×5=4775" "683×9=6147"
Replace-Text "558×5=2790" "521×2=1042"

Replace-Text "324×4=1296" "566×4=2264"
Replace-Text "276×8=2208" "186×4=744"
Replace-Text "448×5=2240" "439×7=3073"
Replace-Text "875×7=6125" "254×5=1270"
Replace-Text "395×9=3555" "989×3=2967"

Replace-Text "411×6=2466" "848×9=7632"
Replace-Text "847×6=5082" "649×9=5841"
Replace-Text "797×4=3188" "705×9=6345"
Replace-Text "148×2=296" "144×9=1296"
Replace-Text "962×2=1924" "918×7=6426"

Write-Host "All replacements complete"
